# "Determinação dos tamanhos de cada armário"
#
# The author re-ordered two shapes on slide 1 ("Cabeamento estruturado"
# floor plan): the background floor-plan picture (Id 3, "Imagem 2") was
# sent to the back of the stack, and the dashed oval annotation
# (Id 177, "Oval 176") used to mark/measure a cabinet area was brought
# back to the front, so the highlight stays visible above every wiring
# element while the floor plan underpins everything else.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

# Floor plan picture -> back of the z-order.
$floorPlan = Get-ShapeById $s 3
if ($floorPlan -ne $null) {
    $floorPlan.ZOrder(1)   # msoSendToBack
}

# Oval used to outline/size the cabinet -> front of the z-order.
$oval = Get-ShapeById $s 177
if ($oval -ne $null) {
    $oval.ZOrder(0)        # msoBringToFront
}
